$wb = $excel.ActiveWorkbook

# New header labels for row 1, columns C..T (A=Depth, B=DRM stay unchanged).
# Values taken from the target (post-edit) state of table/env_corr.xlsx.
$headers = @{
    "C1" = "Temperature"
    "D1" = "Salinity"
    "E1" = "SigmaTheta"
    "F1" = "Density"
    "G1" = "Oxygen"
    "H1" = "Fluorescence"
    "I1" = "Transmission"
    "J1" = "Sand"
    "K1" = "Silt"
    "L1" = "Clay"
    "M1" = "D50"
    "N1" = "TOC"
    "O1" = "TN"
    "P1" = "CN"
    "Q1" = "delta13C"
    "R1" = "Chla"
    "S1" = "WC"
    "T1" = "Porosity"
}

foreach ($ws in $wb.Worksheets) {
    foreach ($addr in $headers.Keys) {
        $ws.Range($addr).Value = $headers[$addr]
    }
}
